$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "55+14=69"
$t.Cell(1,2).Range.Text = "74+20=94"
$t.Cell(1,3).Range.Text = "80+12=92"
$t.Cell(1,4).Range.Text = "92-71=21"
$t.Cell(1,5).Range.Text = "47+8=55"
$t.Cell(2,1).Range.Text = "55-17=38"
$t.Cell(2,2).Range.Text = "96-34=62"
$t.Cell(2,3).Range.Text = "41+34=75"
$t.Cell(2,4).Range.Text = "4+47=51"
$t.Cell(2,5).Range.Text = "34-5=29"
$t.Cell(3,1).Range.Text = "45-0=45"
$t.Cell(3,2).Range.Text = "69-43=26"
$t.Cell(3,3).Range.Text = "28+1=29"
$t.Cell(3,4).Range.Text = "53+18=71"
$t.Cell(3,5).Range.Text = "79-51=28"
$t.Cell(4,1).Range.Text = "68-63=5"
$t.Cell(4,2).Range.Text = "66+9=75"
$t.Cell(4,3).Range.Text = "34-7=27"
$t.Cell(4,4).Range.Text = "39+55=94"
$t.Cell(4,5).Range.Text = "0+55=55"
$t.Cell(5,1).Range.Text = "28+23=51"
$t.Cell(5,2).Range.Text = "42+20=62"
$t.Cell(5,3).Range.Text = "94-32=62"
$t.Cell(5,4).Range.Text = "92-50=42"
$t.Cell(5,5).Range.Text = "10+67=77"
$t.Cell(6,1).Range.Text = "24+41=65"
$t.Cell(6,2).Range.Text = "29+42=71"
$t.Cell(6,3).Range.Text = "61-31=30"
$t.Cell(6,4).Range.Text = "85-58=27"
$t.Cell(6,5).Range.Text = "50-44=6"
$t.Cell(7,1).Range.Text = "13+49=62"
$t.Cell(7,2).Range.Text = "69-1=68"
$t.Cell(7,3).Range.Text = "57-23=34"
$t.Cell(7,4).Range.Text = "15-3=12"
$t.Cell(7,5).Range.Text = "53-22=31"
$t.Cell(8,1).Range.Text = "12+44=56"
$t.Cell(8,2).Range.Text = "85-73=12"
$t.Cell(8,3).Range.Text = "11+69=80"
$t.Cell(8,4).Range.Text = "24+35=59"
$t.Cell(8,5).Range.Text = "58-11=47"
$t.Cell(9,1).Range.Text = "3+40=43"
$t.Cell(9,2).Range.Text = "19-16=3"
$t.Cell(9,3).Range.Text = "77-15=62"
$t.Cell(9,4).Range.Text = "28+25=53"
$t.Cell(9,5).Range.Text = "45-24=21"
$t.Cell(10,1).Range.Text = "1+3=4"
$t.Cell(10,2).Range.Text = "84-31=53"
$t.Cell(10,3).Range.Text = "94-12=82"
$t.Cell(10,4).Range.Text = "84-68=16"
$t.Cell(10,5).Range.Text = "60-58=2"
$t.Cell(11,1).Range.Text = "9+72=81"
$t.Cell(11,2).Range.Text = "59+19=78"
$t.Cell(11,3).Range.Text = "34+33=67"
$t.Cell(11,4).Range.Text = "19+1=20"
$t.Cell(11,5).Range.Text = "95+0=95"
$t.Cell(12,1).Range.Text = "92-31=61"
$t.Cell(12,2).Range.Text = "96-76=20"
$t.Cell(12,3).Range.Text = "82-61=21"
$t.Cell(12,4).Range.Text = "12+19=31"
$t.Cell(12,5).Range.Text = "48+1=49"
$t.Cell(13,1).Range.Text = "60-56=4"
$t.Cell(13,2).Range.Text = "71-21=50"
$t.Cell(13,3).Range.Text = "74+1=75"
$t.Cell(13,4).Range.Text = "81-36=45"
$t.Cell(13,5).Range.Text = "25+23=48"
$t.Cell(14,1).Range.Text = "19+25=44"
$t.Cell(14,2).Range.Text = "91-85=6"
$t.Cell(14,3).Range.Text = "64+14=78"
$t.Cell(14,4).Range.Text = "25+63=88"
$t.Cell(14,5).Range.Text = "26+24=50"
$t.Cell(15,1).Range.Text = "43+32=75"
$t.Cell(15,2).Range.Text = "75-14=61"
$t.Cell(15,3).Range.Text = "53-24=29"
$t.Cell(15,4).Range.Text = "96-24=72"
$t.Cell(15,5).Range.Text = "99-98=1"
$t.Cell(16,1).Range.Text = "21+74=95"
$t.Cell(16,2).Range.Text = "92+6=98"
$t.Cell(16,3).Range.Text = "14+77=91"
$t.Cell(16,4).Range.Text = "46+50=96"
$t.Cell(16,5).Range.Text = "42+57=99"
$t.Cell(17,1).Range.Text = "34+55=89"
$t.Cell(17,2).Range.Text = "16+9=25"
$t.Cell(17,3).Range.Text = "5+62=67"
$t.Cell(17,4).Range.Text = "21+76=97"
$t.Cell(17,5).Range.Text = "48-13=35"
$t.Cell(18,1).Range.Text = "83-49=34"
$t.Cell(18,2).Range.Text = "87-50=37"
$t.Cell(18,3).Range.Text = "39+39=78"
$t.Cell(18,4).Range.Text = "94-80=14"
$t.Cell(18,5).Range.Text = "8-3=5"
$t.Cell(19,1).Range.Text = "92-79=13"
$t.Cell(19,2).Range.Text = "92-29=63"
$t.Cell(19,3).Range.Text = "46-23=23"
$t.Cell(19,4).Range.Text = "29+47=76"
$t.Cell(19,5).Range.Text = "83-69=14"
$t.Cell(20,1).Range.Text = "92-19=73"
$t.Cell(20,2).Range.Text = "38-7=31"
$t.Cell(20,3).Range.Text = "50-39=11"
$t.Cell(20,4).Range.Text = "57-5=52"
$t.Cell(20,5).Range.Text = "49+29=78"
